$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Database")

# Insert a new row at position 3 (pushes existing rows 3..21 down to 4..22)
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with the RESOURCELINK data source entry.
# Set the H3 note before A3 so the shared-string table gains
# "Attempt to standardise..." (191) before "RESOURCELINK" (192).
$ws.Cells.Item(3, 8).Value = "Attempt to standardise the DSN across environments."
$ws.Cells.Item(3, 1).Value = "RESOURCELINK"
$ws.Cells.Item(3, 2).Value = "CMTEST"
$ws.Cells.Item(3, 3).Value = "jm08_cmt"
$ws.Cells.Item(3, 4).Formula = '=CONCATENATE( "mqsisetdbparms ",ConfigData!$D$4," -n ",A3," -u ",B3," -p ",C3)'
$ws.Cells.Item(3, 6).Formula = '=CONCATENATE( "mqsicvp ",ConfigData!$D$4," -n ",A3)'

# Update the view: scroll back to top-left (removes topLeftCell="A7")
# and move the selection to F3 (was F21).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F3").Select()
